$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(821).ClearFormats()
